{"js": "// Diff summary:\n//  - remove the paragraph \"Gestion des contraintes de validations de station ...\"\n//  - add four new TODO paragraphs after \"Enlever la possibilit\u00e9 de cr\u00e9er des capteurs\"\n//    (and before the trailing empty/_GoBack paragraph)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Remove the paragraph about station Nom/Client validation constraints.\nconst obsolete = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Gestion des contraintes de validations de station\") !== -1\n);\nif (obsolete) {\n  obsolete.delete();\n  await context.sync();\n}\n\n// 2) Find the paragraph that should receive the four new items right after it.\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items/text\");\nawait context.sync();\n\nconst anchor = refreshed.items.find((p) =>\n  p.text.indexOf(\"Enlever la possibilit\u00e9 de cr\u00e9er des capteurs\") !== -1\n);\nconst insertionRange = anchor.getRange(Word.RangeLocation.end);\n\n// 3) Insert the four new paragraphs as one OOXML blob so the runs/proofErr markers\n//    (spell-check wrappers around \"checkboxes\", \"acceuil\", \"acces\", \"nbutilisateur\",\n//    \"parametre\") come out exactly as Word itself would produce them.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p><w:r><w:t>Mettre ordre alphab\\u00e9tique tous les menus</w:t></w:r></w:p>' +\n            '<w:p>' +\n              '<w:r><w:t xml:space=\"preserve\">Gestion des </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>checkboxes</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\"> de suppression</w:t></w:r>' +\n            '</w:p>' +\n            '<w:p>' +\n              '<w:r><w:t>Page d\\u2019</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>acceuil</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\"> en libre </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>acces</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\"> + stations publics</w:t></w:r>' +\n            '</w:p>' +\n            '<w:p>' +\n              '<w:r><w:t xml:space=\"preserve\">Admin menu =&gt; </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>nbutilisateur</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\"> / limitation </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>parametre</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ninsertionRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the paragraph about station Nom/Client validation constraints.\n$findRange = $d.Content\nif ($findRange.Find.Execute(\"Gestion des contraintes de validations de station\")) {\n    $obsoleteParagraph = $findRange.Paragraphs(1)\n    $obsoleteParagraph.Range.Delete()\n}\n\n# 2) Locate the paragraph \"Enlever la possibilit\u00e9 de cr\u00e9er des capteurs\" \u2014 the four new\n#    TODO paragraphs are inserted right after it, before the trailing empty/_GoBack\n#    paragraph.\n$anchorRange = $d.Content\n$null = $anchorRange.Find.Execute(\"Enlever la possibilit\u00e9 de cr\u00e9er des capteurs\")\n$anchorParagraph = $anchorRange.Paragraphs(1)\n$insertionPoint = $d.Range($anchorParagraph.Range.End, $anchorParagraph.Range.End)\n\n# 3) Insert the four new paragraphs as one WordOpenXML blob so the runs/proofErr\n#    markers (spell-check wrappers around \"checkboxes\", \"acceuil\", \"acces\",\n#    \"nbutilisateur\", \"parametre\") come out exactly as Word itself would produce them.\n#    InsertXML needs a trailing empty <w:p/> to supply the paragraph mark that splits\n#    the new content away from the following paragraph; that spare paragraph is\n#    removed again afterwards.\n$xml = @'\n<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p><w:r><w:t>Mettre ordre alphab\u00e9tique tous les menus</w:t></w:r></w:p>\n<w:p><w:r><w:t xml:space=\"preserve\">Gestion des </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>checkboxes</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> de suppression</w:t></w:r></w:p>\n<w:p><w:r><w:t>Page d\u2019</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>acceuil</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> en libre </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>acces</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> + stations publics</w:t></w:r></w:p>\n<w:p><w:r><w:t xml:space=\"preserve\">Admin menu =&gt; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>nbutilisateur</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> / limitation </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>parametre</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>\n<w:p/>\n</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n$insertionPoint.InsertXML($xml)\n\n$extraIndex = $d.Paragraphs.Count - 1\n$d.Paragraphs($extraIndex).Range.Delete()\n"}
